# Cambios Generales Junio 2024
$wb = $excel.ActiveWorkbook

$wsPlantilla  = $wb.Worksheets.Item("Plantilla")
$wsSucursales = $wb.Worksheets.Item("Sucursales")

# --- Plantilla: add new column V with header "Incluye Producto(1:SI,0:No)" ---
$wsPlantilla.Range("V1").Value = "Incluye Producto(1:SI,0:No)"

# Match the formatting of the existing header row (e.g. Q1) for the new header cell
$wsPlantilla.Range("Q1").Copy() | Out-Null
$wsPlantilla.Range("V1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Give the new column a sensible width like the rest of the layout
# (target stored width is 24.28515625 chars / 170px; engine quantizes
# ColumnWidth assignment to pixel steps, so feed it the char-width value
# whose pixel-rounded result lands closest to that target)
$wsPlantilla.Range("V1").ColumnWidth = 23.45

# --- Selections / active sheet as left by the editing session ---
$wsPlantilla.Range("U2").Select() | Out-Null

$wsSucursales.Activate()
$wsSucursales.Range("A2").Select() | Out-Null
